$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (weekly price update), shifting the
# existing rows 59-81 down to 60-82.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly price record.
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44704
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = 100114002
$ws.Range("G59").Value = "Camote"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 80
$ws.Range("K59").Value = 18000
$ws.Range("L59").Value = 18000
$ws.Range("M59").Value = 18000
$ws.Range("N59").Value = "$/malla 20 kilos"
$ws.Range("O59").Value = "Perú"
$ws.Range("P59").Value = 900
$ws.Range("Q59").Value = 20
$ws.Range("R59").Value = "Hortaliza"
